$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column Q (17): old Q,R,S shift to T,U,V
$ws.Range("Q1:S1").EntireColumn.Insert()

# New headers for the inserted "Grupa (prowadzacy)" columns
$ws.Range("Q1").Value = "Grupa (prowadzący):"
$ws.Range("R1").Value = "Punkty — Grupa (prowadzący):"
$ws.Range("S1").Value = "Opinia — Grupa (prowadzący):"

# New instructor/group values for the three response rows
$ws.Range("Q2").Value = "Zbigniew Kaleta"
$ws.Range("Q3").Value = "Michał Idzik"
$ws.Range("Q4").Value = "Bernard Maj"

# Re-point the active selection / view like the target file
$ws.Range("R7").Select()

# ---- formatting pass -------------------------------------------------
# Plain "label" style used by headers and most text cells: Calibri 11, black RGB
$labelRange = $ws.Range("A1:V1,H2:H4,K2:K4,N2:N4,T2:T4,U2:U4,V2:V4")
$labelRange.Font.Name = "Calibri"
$labelRange.Font.Size = 11
$labelRange.Font.Color = 0x000000

# "Grupa" column style: Arial 10, left aligned
$grupaRange = $ws.Range("Q1:S1,Q2:Q4")
$grupaRange.Font.Name = "Arial"
$grupaRange.Font.Size = 10
$grupaRange.HorizontalAlignment = -4131

# Numeric value cells: Calibri 11 black, right aligned
$numRange = $ws.Range("A2:A4,F2:F4,I2:I4,L2:L4")
$numRange.Font.Name = "Calibri"
$numRange.Font.Size = 11
$numRange.Font.Color = 0x000000
$numRange.HorizontalAlignment = -4152

# Date value cells: Calibri 11 black, right aligned, keep existing date number format
$dateRange = $ws.Range("B2:C4")
$dateRange.Font.Name = "Calibri"
$dateRange.Font.Size = 11
$dateRange.Font.Color = 0x000000
$dateRange.HorizontalAlignment = -4152

# Email value cells: underlined hyperlink-blue Calibri 11
$emailRange = $ws.Range("D2:E4")
$emailRange.Font.Name = "Calibri"
$emailRange.Font.Size = 11
$emailRange.Font.Underline = 2
$emailRange.Font.Color = 0x0563C1

Write-Output "done"
